# "Fruta / hortaliza, semanal" weekly update.
# A new weekly price record is inserted at row 199 (pushing the existing
# rows 199-271 down to 200-272); the new record copies the category/
# quality/unit template of the most recent "Segunda" entry but carries a
# new (later) observation date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 199; everything below
# (old rows 199-271) shifts down by one, landing on rows 200-272.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with this week's data point.
$ws.Cells.Item(199, 1).Value  = 7
$ws.Cells.Item(199, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(199, 3).Value  = 'Ñuble'
$ws.Cells.Item(199, 4).Value  = 44924
$ws.Cells.Item(199, 5).Value  = 16
$ws.Cells.Item(199, 6).Value  = 'Fruta'
$ws.Cells.Item(199, 7).Value  = 100108
$ws.Cells.Item(199, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(199, 9).Value  = 100108005
$ws.Cells.Item(199, 10).Value = 'Piña'
$ws.Cells.Item(199, 11).Value = 'Caramelo'
$ws.Cells.Item(199, 12).Value = 'Segunda'
$ws.Cells.Item(199, 13).Value = 120
$ws.Cells.Item(199, 14).Value = 19000
$ws.Cells.Item(199, 15).Value = 20000
$ws.Cells.Item(199, 16).Value = 19500
$ws.Cells.Item(199, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(199, 18).Value = 'Ecuador'
$ws.Cells.Item(199, 19).Value = 1393
$ws.Cells.Item(199, 20).Value = 14
